# Added controllers for randomMoveSymbols
# The symbol weight table (rows 2-25, columns A:F) has been shuffled: several
# rows' complete 6-value records have been moved to different row positions.
# Apply the new values cell-by-cell so the resulting sheet matches the target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $a, $b, $c, $d, $e, $f) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
}

Set-Row 4  1001 18 30 75 60 72
Set-Row 5  501  9  52 30 75 45
Set-Row 6  701  3  90 45 97 15
Set-Row 7  901  16 15 45 60 60
Set-Row 8  301  6  45 30 60 45
Set-Row 9  401  9  48 67 75 45
Set-Row 10 1201 2  10 10 10 10
Set-Row 11 1203 3  15 15 15 15

Set-Row 14 201  9  30 15 45 30
Set-Row 15 1202 2  10 10 10 10
Set-Row 16 802  0  4  5  4  0
Set-Row 17 3    0  3  3  3  3
Set-Row 18 1    0  2  2  2  2
Set-Row 19 1101 0  15 30 30 0
Set-Row 20 2    0  2  2  2  2
Set-Row 21 502  0  4  0  0  0
Set-Row 22 402  0  0  4  0  0
Set-Row 23 602  0  0  4  0  9
